$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Custom Table Entry")

$ws.Range("A4").Value = "4b44pz84e13934k"
$ws.Range("A5").Value = "15ikn58z719910o"
$ws.Range("A6").Value = "1b1ejz1ds8x9782"
$ws.Range("A7").Value = "ir6w36l8441kx3f"

$ws.Range("A9").Value = "4a73q62w180m377"
$ws.Range("A10").Value = "ut3j26acgfkqnuo"
$ws.Range("A11").Value = "061atd70d46krvg"

$ws.Range("A13").Value = "p9x3i35263e933x"
